# "Fruta / hortaliza, semanal" weekly refresh:
# A new weekly price-report row is inserted at the top of this product's
# data block (row 82), pushing the existing rows 82-165 down to 83-166
# (dimension grows from A1:R165 to A1:R166). The brand-new row duplicates
# the values that used to sit in row 82 (same volume/price/quality figures)
# but is dated one day after the most recent existing record (44704 -> 44705),
# i.e. it is simply this week's placeholder entry stamped with the new date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the first data row of this block; this shifts
# everything that was in rows 82:165 down into rows 83:166 automatically.
$ws.Rows.Item(82).Insert()

# The row that is now at 83 holds exactly what used to be row 82 - copy it
# back up into the freshly-inserted row 82 so the new entry starts out as a
# duplicate of last week's figures.
$ws.Range("A83:R83").Copy($ws.Range("A82"))

# ...then stamp the new entry with this week's date (one day after the
# latest date already present in the sheet, 44704 -> 44705).
$ws.Range("D82").Value = 44705
